$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.360.82"
$ws.Range("E2").Value = "  -4.62%  "
$ws.Range("D3").Value = "1.564.78"
$ws.Range("E3").Value = "  -4.98%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3707"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3389"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07646"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.28%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.047"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.925"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.79%  "
$ws.Range("D16").Value = "1.568.49"
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001128"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06734"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.244"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.83%  "
$ws.Range("E22").Value = "  -5.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5275"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.99%  "
$ws.Range("D25").Value = "22.357.91"
$ws.Range("E25").Value = "  -4.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.361"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.821"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.983"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("D32").Value = "1.732.49"
$ws.Range("E32").Value = "  -5.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.230"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.27%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.007"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.006"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08485"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02534"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2322"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.534"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06413"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.289"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -8.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6335"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5980"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.89%  "
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.099"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.266"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.29%  "
